# Auto-generated edit script: updates crypto price/volume table cells
# to match refreshed data feed values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.452.56"
$ws.Range("E2").Value = "'  +0.30%  "
$ws.Range("D3").Value = "'1.806.12"
$ws.Range("E3").Value = "'  +0.09%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'225.49"
$ws.Range("E5").Value = "'  -0.84%  "
$ws.Range("D6").Value = "'0.587"
$ws.Range("E6").Value = "'  +2.10%  "
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'38.19"
$ws.Range("E8").Value = "'  +5.85%  "
$ws.Range("D9").Value = "'0.287"
$ws.Range("E9").Value = "'  -4.76%  "
$ws.Range("E10").Value = "'  -3.37%  "
$ws.Range("E11").Value = "'  +0.97%  "
$ws.Range("D12").Value = "'2.068.41"
$ws.Range("E12").Value = "'  +0.15%  "
$ws.Range("D13").Value = "'11.08"
$ws.Range("E13").Value = "'  -5.70%  "
$ws.Range("D14").Value = "'1.812.78"
$ws.Range("E14").Value = "'  +0.40%  "
$ws.Range("D15").Value = "'34.435.35"
$ws.Range("E15").Value = "'  +0.29%  "
$ws.Range("D16").Value = "'0.627"
$ws.Range("E16").Value = "'  -2.88%  "
$ws.Range("D17").Value = "'4.39"
$ws.Range("E17").Value = "'  -2.74%  "
$ws.Range("D18").Value = "'67.87"
$ws.Range("E18").Value = "'  -1.79%  "
$ws.Range("D19").Value = "'241.73"
$ws.Range("E19").Value = "'  -1.55%  "
$ws.Range("D20").Value = "'0.0₃0766"
$ws.Range("E20").Value = "'  -3.63%  "
$ws.Range("D21").Value = "'11.07"
$ws.Range("E21").Value = "'  -4.47%  "
$ws.Range("E22").Value = "'  +0.04%  "
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "'  -2.38%  "
$ws.Range("E24").Value = "'  +3.65%  "
$ws.Range("D25").Value = "'170.16"
$ws.Range("E25").Value = "'  -0.97%  "
$ws.Range("D26").Value = "'7.70"
$ws.Range("E26").Value = "'  -3.73%  "
$ws.Range("D27").Value = "'17.42"
$ws.Range("E27").Value = "'  +3.08%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "'  +0.84%  "
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("E30").Value = "'  -1.04%  "
$ws.Range("D31").Value = "'3.75"
$ws.Range("E31").Value = "'  -2.56%  "
$ws.Range("D32").Value = "'0.0512"
$ws.Range("E32").Value = "'  -3.74%  "
$ws.Range("D33").Value = "'3.81"
$ws.Range("E33").Value = "'  -5.32%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "'  -0.97%  "
$ws.Range("D35").Value = "'1.328.86"
$ws.Range("E35").Value = "'  -4.47%  "
$ws.Range("D36").Value = "'0.638"
$ws.Range("E36").Value = "'  -5.26%  "
$ws.Range("E37").Value = "'  -1.17%  "
$ws.Range("D38").Value = "'0.0188"
$ws.Range("E38").Value = "'  -1.31%  "
$ws.Range("D39").Value = "'2.31"
$ws.Range("E39").Value = "'  -6.18%  "
$ws.Range("E40").Value = "'  +1.47%  "
$ws.Range("E41").Value = "'  -1.47%  "
$ws.Range("D42").Value = "'81.57"
$ws.Range("E42").Value = "'  -1.30%  "
$ws.Range("B43").Value = "'MXToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.79"
$ws.Range("E43").Value = "'  -1.15%  "
$ws.Range("B44").Value = "'ARBITRUM"
$ws.Range("C44").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.943"
$ws.Range("E44").Value = "'  -2.39%  "
$ws.Range("D45").Value = "'13.53"
$ws.Range("E45").Value = "'  +0.66%  "
$ws.Range("E46").Value = "'  +1.53%  "
$ws.Range("D47").Value = "'1.968.55"
$ws.Range("E47").Value = "'  +0.10%  "
$ws.Range("D48").Value = "'5.71"
$ws.Range("E48").Value = "'  -5.01%  "
$ws.Range("E49").Value = "'  +0.05%  "
$ws.Range("D50").Value = "'101.82"
$ws.Range("E50").Value = "'  -2.62%  "
$ws.Range("E51").Value = "'  -5.33%  "
